$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Candidate ID (row 2)
$ws.Range("B2").Value = 23081721

# Update Client Id (row 2)
$ws.Range("A2").Value = "MXUWy809"

# Update User Name (row 2)
$ws.Range("C2").Value = "lgfhjxc22"

# Update Exam Password (row 2)
$ws.Range("D2").Value = "cs!5&TB8"

# Update First Name (row 2)
$ws.Range("F2").Value = "MhwkBWYz"

# Update Last Name (row 2)
$ws.Range("G2").Value = "zeWE"
